function Set-NumericText {
    param($ws, $cellAddr, $val)
    # Force the value to be stored as text even though it looks numeric,
    # then clear the resulting cell format so no stray style is left behind.
    $c = $ws.Range($cellAddr)
    $c.Value = "'" + $val
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.391.36'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.841.64'
$ws.Range("E3").Value = '  -0.30%  '
Set-NumericText $ws "D4" '0.9989'
$ws.Range("E4").Value = '  +0.10%  '
Set-NumericText $ws "D5" '239.60'
$ws.Range("E5").Value = '  -0.31%  '
Set-NumericText $ws "D6" '0.6265'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -0.76%  '
Set-NumericText $ws "D9" '0.2900'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("E10").Value = '  +1.34%  '
Set-NumericText $ws "D11" '0.07722'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '1.830.11'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("E15").Value = '  -3.08%  '
Set-NumericText $ws "D16" '81.75'
Set-NumericText $ws "D17" '6.243'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '29.437.28'
$ws.Range("E18").Value = '  +0.02%  '
Set-NumericText $ws "D19" '231.63'
$ws.Range("E19").Value = '  +0.92%  '
Set-NumericText $ws "D20" '12.31'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("E23").Value = '  +0.08%  '
Set-NumericText $ws "D24" '158.40'
$ws.Range("E24").Value = '  -0.48%  '
Set-NumericText $ws "D25" '8.481'
$ws.Range("E25").Value = '  +0.84%  '
Set-NumericText $ws "D26" '0.1351'
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("E27").Value = '  -0.96%  '
Set-NumericText $ws "D28" '0.07323'
$ws.Range("E28").Value = '  +13.64%  '
Set-NumericText $ws "D29" '1.465'
$ws.Range("E29").Value = '  +3.72%  '
Set-NumericText $ws "D30" '1.478'
$ws.Range("E30").Value = '  +0.21%  '
# Row 31/32: coins swapped order (Filecoin now ranks above InternetComputer)
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-NumericText $ws "D31" '4.051'
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-NumericText $ws "D32" '4.057'
$ws.Range("E32").Value = '  -0.97%  '

Set-NumericText $ws "D33" '1.815'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -0.19%  '
Set-NumericText $ws "D35" '0.6974'
$ws.Range("E35").Value = '  +0.13%  '
Set-NumericText $ws "D36" '2.566'
$ws.Range("E36").Value = '  -0.49%  '
Set-NumericText $ws "D37" '6.962'
$ws.Range("E37").Value = '  +3.60%  '
Set-NumericText $ws "D38" '0.01838'
$ws.Range("E38").Value = '  +0.14%  '
Set-NumericText $ws "D39" '2.815'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '1.235.66'
$ws.Range("E40").Value = '  -2.75%  '
Set-NumericText $ws "D41" '0.9454'
$ws.Range("E41").Value = '  +4.01%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").Value = '1.996.02'
$ws.Range("E43").Value = '  -0.67%  '
Set-NumericText $ws "D44" '100.74'
$ws.Range("E44").Value = '  -0.56%  '
Set-NumericText $ws "D45" '65.65'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("E46").Value = '  +0.24%  '
Set-NumericText $ws "D47" '1.722'
$ws.Range("E47").Value = '  -1.18%  '
Set-NumericText $ws "D48" '6.959'
$ws.Range("E48").Value = '  -1.68%  '
Set-NumericText $ws "D49" '8.959'
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("E51").Value = '  -1.41%  '
